$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 4 (DownloadEmailAttachments) value from "Yes" to "No"
$ws.Range("B4").Value = "No"

# Add new row 7 for FormatWordReport setting
$ws.Range("A7").Value = "FormatWordReport"
$ws.Range("B7").Value = "No"
$ws.Range("C7").Value = "Yes - Format paragraphs into table"

# Update the active selection to match the post-edit state
$ws.Range("B8").Select()
